$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume(1h) data, one @{Row=; D=; E=} entry per changed row.
$updates = @(
    @{ Row = 2; D = '62.940.21'; E = '  -0.09%  ' }
    @{ Row = 3; D = '2.583.25'; E = '  +1.24%  ' }
    @{ Row = 4; D = $null; E = '  +0.01%  ' }
    @{ Row = 5; D = '582.85'; E = '  +0.50%  ' }
    @{ Row = 6; D = '146.76'; E = '  -0.24%  ' }
    @{ Row = 7; D = $null; E = '  +0.01%  ' }
    @{ Row = 8; D = $null; E = '  +2.12%  ' }
    @{ Row = 9; D = $null; E = '  +2.40%  ' }
    @{ Row = 10; D = $null; E = '  +2.55%  ' }
    @{ Row = 11; D = $null; E = '  -0.13%  ' }
    @{ Row = 12; D = $null; E = '  -0.22%  ' }
    @{ Row = 13; D = '27.34'; E = '  +0.49%  ' }
    @{ Row = 14; D = '3.045.77'; E = '  +1.22%  ' }
    @{ Row = 15; D = '62.814.05'; E = '  -0.15%  ' }
    @{ Row = 16; D = $null; E = '  +3.09%  ' }
    @{ Row = 17; D = '2.596.26'; E = '  +1.35%  ' }
    @{ Row = 18; D = '11.30'; E = '  -0.27%  ' }
    @{ Row = 19; D = '341.95'; E = '  +1.87%  ' }
    @{ Row = 20; D = $null; E = '  +0.85%  ' }
    @{ Row = 21; D = $null; E = '  -0.79%  ' }
    @{ Row = 22; D = $null; E = '  +0.02%  ' }
    @{ Row = 23; D = '67.07'; E = '  +2.35%  ' }
    @{ Row = 24; D = '2.708.67'; E = '  +1.15%  ' }
    @{ Row = 25; D = $null; E = '  -1.51%  ' }
    @{ Row = 26; D = $null; E = '  -1.22%  ' }
    @{ Row = 27; D = '1.00'; E = '  +0.06%  ' }
    @{ Row = 28; D = $null; E = '  -0.68%  ' }
    @{ Row = 29; D = '7.84'; E = '  +7.02%  ' }
    @{ Row = 30; D = $null; E = '  -2.34%  ' }
    @{ Row = 31; D = '1.92'; E = '  +1.63%  ' }
    @{ Row = 32; D = '0.0₃0821'; E = '  +1.07%  ' }
    @{ Row = 33; D = '467.89'; E = '  +13.87%  ' }
    @{ Row = 34; D = '175.60'; E = '  -0.90%  ' }
    @{ Row = 35; D = '1.60'; E = '  +3.72%  ' }
    @{ Row = 36; D = $null; E = '  +0.11%  ' }
    @{ Row = 37; D = $null; E = '  +0.30%  ' }
    @{ Row = 38; D = $null; E = '  -0.65%  ' }
    @{ Row = 39; D = $null; E = '  +4.25%  ' }
    @{ Row = 41; D = $null; E = '  -1.91%  ' }
    @{ Row = 42; D = '157.57'; E = '  +4.23%  ' }
    @{ Row = 43; D = $null; E = '  +0.06%  ' }
    @{ Row = 44; D = '0.633'; E = '  +5.00%  ' }
    @{ Row = 45; D = $null; E = '  +1.57%  ' }
    @{ Row = 46; D = '0.0541'; E = '  +0.57%  ' }
    @{ Row = 47; D = '0.0966'; E = '  -0.07%  ' }
    @{ Row = 48; D = $null; E = '  -0.86%  ' }
    @{ Row = 49; D = '18.34'; E = '  +0.52%  ' }
    @{ Row = 50; D = $null; E = '  +0.39%  ' }
    @{ Row = 51; D = '11.42'; E = '  +0.98%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
